$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 16 - this shifts the existing rows 16..200 down to 17..201
# (old row 16 -> new row 17, old row 17 -> new row 18, ..., old row 200 (the
# thick-bottom totals row) -> new row 201), and grows the running-sum/shared
# "=D-C" formula ranges to cover it.
$ws.Rows(16).Insert()

# --- Row 16: brand-new "Research" entry -----------------------------------
# Pick up the date/time number formats from row 15 (a fully filled-in row)
# instead of leaving the generic "blank template" formatting behind.
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)

$ws.Range("A16").Value = 45969
$ws.Range("B16").Value = "Research"
$ws.Range("C16").Value = 0.77222222222222225
$ws.Range("D16").Value = 0.86458333333333337
$ws.Range("E16").Formula = "=D16-C16"

# --- Row 17: old (blank template) row 16, shifted down - completed -------
# "Check-in" entry. Only the date/start/end time columns need their number
# format fixed up; category and description columns already use the right
# style.
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("C15:D15").Copy()
$ws.Range("C17:D17").PasteSpecial(-4122)

$ws.Range("A17").Value = 45971
$ws.Range("B17").Value = "Check-in"
$ws.Range("C17").Value = 0.625
$ws.Range("D17").Value = 0.63888888888888884
# Set this description before row 16's, so the new shared-string entries end
# up in the same order the workbook author originally typed them in.
$ws.Range("F17").Value = "Check-in: Group meeting with Tyler, established communicaton between frontend and backend and discussed tasks to complete for the prototype"

$ws.Range("F16").Value = "Researching communication from backend to frontend and Flutter mobile app components. GET api call works and can properly pull unformatted data"

# --- Row 18: old (blank template) row 17, shifted down - still in-progress
# "Frontend Dev" entry: only date, category and start time are filled in so
# far, end time/description are left blank.
$ws.Range("A15").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("A18").Value = 45972
$ws.Range("B18").Value = "Frontend Dev"
$ws.Range("C18").Value = 0.45902777777777776

# Restore the view: scrolled down a bit, with F18 selected (the still-open
# Description cell for the in-progress entry).
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("F18").Select()
